# Update the date header and the 25 "three-digit ÷ one-digit" practice
# answers in the table to the new day's values.
#
# Every "FindText" value below is unique within the document, so a plain
# (non-wildcard) Find/Replace targeted at the whole document content is
# sufficient and safe - it cannot clobber an unrelated cell.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-05-31 Saturday"; New = "2025-06-01 Sunday" },
    @{ Old = "656÷3=218, 2";        New = "874÷9=97, 1" },
    @{ Old = "208÷5=41, 3";         New = "342÷9=38, 0" },
    @{ Old = "658÷3=219, 1";        New = "438÷4=109, 2" },
    @{ Old = "191÷7=27, 2";         New = "268÷3=89, 1" },
    @{ Old = "757÷2=378, 1";        New = "565÷4=141, 1" },
    @{ Old = "626÷4=156, 2";        New = "624÷9=69, 3" },
    @{ Old = "230÷8=28, 6";         New = "593÷2=296, 1" },
    @{ Old = "242÷3=80, 2";         New = "703÷6=117, 1" },
    @{ Old = "738÷3=246, 0";        New = "593÷6=98, 5" },
    @{ Old = "523÷6=87, 1";         New = "997÷6=166, 1" },
    @{ Old = "483÷5=96, 3";         New = "226÷3=75, 1" },
    @{ Old = "152÷4=38, 0";         New = "763÷9=84, 7" },
    @{ Old = "943÷6=157, 1";        New = "658÷5=131, 3" },
    @{ Old = "892÷2=446, 0";        New = "454÷8=56, 6" },
    @{ Old = "958÷5=191, 3";        New = "118÷3=39, 1" },
    @{ Old = "871÷5=174, 1";        New = "521÷4=130, 1" },
    @{ Old = "216÷3=72, 0";         New = "664÷4=166, 0" },
    @{ Old = "659÷4=164, 3";        New = "973÷4=243, 1" },
    @{ Old = "950÷4=237, 2";        New = "939÷3=313, 0" },
    @{ Old = "485÷3=161, 2";        New = "957÷9=106, 3" },
    @{ Old = "569÷2=284, 1";        New = "696÷9=77, 3" },
    @{ Old = "786÷2=393, 0";        New = "654÷6=109, 0" },
    @{ Old = "134÷5=26, 4";         New = "769÷6=128, 1" },
    @{ Old = "151÷5=30, 1";         New = "921÷9=102, 3" },
    @{ Old = "587÷6=97, 5";         New = "512÷9=56, 8" }
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $pair.New, 2) | Out-Null
}
